# Auto-generated edit script applying the profit-table refresh diff
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 2865.9092
$ws.Range("I13").Value = 1525
$ws.Range("K13").Value = 1525
$ws.Range("M13").Value = -1356
$ws.Range("H18").Value = 647
$ws.Range("I18").Value = 602.875
$ws.Range("K18").Value = 602.875
$ws.Range("M18").Value = -318.875
$ws.Range("H100").Value = 1565.7142
$ws.Range("I100").Value = 1542.5
$ws.Range("J100").Value = 1640
$ws.Range("K100").Value = 1542.5
$ws.Range("L100").Value = 1640
$ws.Range("M100").Value = -1001.5
$ws.Range("N100").Value = -2722
$ws.Range("J112").Value = 71430160
$ws.Range("L112").Value = 214290480
$ws.Range("N112").Value = -214292696
$ws.Range("H129").Value = 1070.0233
$ws.Range("J129").Value = 1170.7297
$ws.Range("L129").Value = 3512.189100000001
$ws.Range("N129").Value = -13512.1891
$ws.Range("H133").Value = 59500
$ws.Range("J133").Value = 59500
$ws.Range("L133").Value = 59500
$ws.Range("N133").Value = -69620
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()  # was -50980
$ws.Range("H137").Value = 66671252
$ws.Range("I137").Value = 2254.7273
$ws.Range("J137").Value = 250011000
$ws.Range("K137").Value = 6764.1819
$ws.Range("L137").Value = 750033000
$ws.Range("M137").Value = -4214.1819
$ws.Range("N137").Value = -750038100
$ws.Range("H138").Value = 5005.8335
$ws.Range("I138").Value = 4080
$ws.Range("J138").Value = 5387.0586
$ws.Range("K138").Value = 12240
$ws.Range("L138").Value = 16161.1758
$ws.Range("M138").Value = -7100
$ws.Range("N138").Value = -26441.1758
# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23444.572
$ws.Range("I32").Value = 18036.46
$ws.Range("K32").Value = 18036.46
$ws.Range("M32").Value = -17749.46
$ws.Range("H37").Value = 10326.625
$ws.Range("J37").Value = 10326.625
$ws.Range("L37").Value = 10326.625
$ws.Range("N37").Value = -10872.625
$ws.Range("H61").Value = 5850774.5
$ws.Range("I61").Value = 7409140.5
$ws.Range("J61").Value = 6901
$ws.Range("K61").Value = 7409140.5
$ws.Range("L61").Value = 6901
$ws.Range("M61").Value = -7408928.5
$ws.Range("N61").Value = -7325
$ws.Range("H132").Value = 1439441.1
$ws.Range("I132").Value = 1569936
$ws.Range("J132").Value = 3997
$ws.Range("K132").Value = 4709808
$ws.Range("L132").Value = 11991
$ws.Range("M132").Value = -4707278
$ws.Range("N132").Value = -17051
$ws.Range("H136").Value = 5850774.5
$ws.Range("I136").Value = 7409140.5
$ws.Range("J136").Value = 6901
$ws.Range("K136").Value = 22227421.5
$ws.Range("L136").Value = 20703
$ws.Range("M136").Value = -22224871.5
$ws.Range("N136").Value = -25803
# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14605.06
$ws.Range("I31").Value = 21234.32
$ws.Range("J31").Value = 7975.8
$ws.Range("K31").Value = 21234.32
$ws.Range("L31").Value = 7975.8
$ws.Range("M31").Value = -20939.32
$ws.Range("N31").Value = -8565.799999999999
$ws.Range("H34").Value = 14605.06
$ws.Range("I34").Value = 21234.32
$ws.Range("J34").Value = 7975.8
$ws.Range("K34").Value = 21234.32
$ws.Range("L34").Value = 7975.8
$ws.Range("M34").Value = -21032.32
$ws.Range("N34").Value = -8379.799999999999
$ws.Range("H41").Value = 5658.3335
$ws.Range("J41").Value = 8666.666999999999
$ws.Range("L41").Value = 8666.666999999999
$ws.Range("N41").Value = -9522.666999999999
$ws.Range("H50").Value = 10359.375
$ws.Range("J50").Value = 10359.375
$ws.Range("L50").Value = 10359.375
$ws.Range("N50").Value = -11609.375
$ws.Range("H51").Value = 10097.6
$ws.Range("J51").Value = 10422
$ws.Range("L51").Value = 10422
$ws.Range("N51").Value = -11894
$ws.Range("H53").Value = 30000
$ws.Range("J53").Value = 30000
$ws.Range("L53").Value = 30000
$ws.Range("N53").Value = -31214
$ws.Range("H59").Value = 14680
$ws.Range("J59").Value = 15608.889
$ws.Range("L59").Value = 15608.889
$ws.Range("N59").Value = -17898.889
$ws.Range("H60").Value = 9236.817999999999
$ws.Range("I60").Value = 5750
$ws.Range("J60").Value = 10011.667
$ws.Range("K60").Value = 5750
$ws.Range("L60").Value = 10011.667
$ws.Range("M60").Value = -5239
$ws.Range("N60").Value = -11033.667
$ws.Range("H61").Value = 10097.6
$ws.Range("J61").Value = 10422
$ws.Range("L61").Value = 10422
$ws.Range("N61").Value = -11118
$ws.Range("H68").Value = 18384.182
$ws.Range("J68").Value = 18795.8
$ws.Range("L68").Value = 18795.8
$ws.Range("N68").Value = -20293.8
$ws.Range("H71").Value = 18384.182
$ws.Range("J71").Value = 18795.8
$ws.Range("L71").Value = 56387.39999999999
$ws.Range("N71").Value = -63875.39999999999
$ws.Range("H74").Value = 15160.818
$ws.Range("J74").Value = 16455.445
$ws.Range("L74").Value = 16455.445
$ws.Range("N74").Value = -18203.445
$ws.Range("H77").Value = 15160.818
$ws.Range("J77").Value = 16455.445
$ws.Range("L77").Value = 49366.335
$ws.Range("N77").Value = -58102.335
$ws.Range("H111").Value = 49702
$ws.Range("J111").Value = 49702
$ws.Range("L111").Value = 49702
$ws.Range("N111").Value = -57882
# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 178402190
$ws.Range("J105").Value = 223002240
$ws.Range("L105").Value = 669006720
$ws.Range("N105").Value = -669011962
$ws.Range("H113").Value = 636.5833
$ws.Range("I113").Value = 709.75
$ws.Range("J113").Value = 600
$ws.Range("K113").Value = 2129.25
$ws.Range("L113").Value = 1800
$ws.Range("M113").Value = 40.75
$ws.Range("N113").Value = -6140
$ws.Range("H122").Value = 5581464.5
$ws.Range("J122").Value = 745749
$ws.Range("L122").Value = 6711741
$ws.Range("N122").Value = -6716641
$ws.Range("H127").Value = 1410.25
$ws.Range("J127").Value = 1410.25
$ws.Range("L127").Value = 4230.75
$ws.Range("N127").Value = -14150.75
$ws.Range("H131").Value = 46671664
$ws.Range("I131").Value = 95246860
$ws.Range("J131").Value = 27781306
$ws.Range("K131").Value = 285740580
$ws.Range("L131").Value = 83343918
$ws.Range("M131").Value = -285735540
$ws.Range("N131").Value = -83353998
$ws.Range("H139").Value = 1931.4286
$ws.Range("I139").Value = 1931.4286
$ws.Range("K139").Value = 5794.2858
$ws.Range("M139").Value = -654.2857999999997
# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 59800
$ws.Range("J133").Value = 59800
$ws.Range("L133").Value = 59800
$ws.Range("N133").Value = -69920
$ws.Range("H140").Value = 57000
$ws.Range("J140").Value = 57000
$ws.Range("L140").Value = 57000
$ws.Range("N140").Value = -67360
# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4701
$ws.Range("I136").Value = 2268
$ws.Range("J136").Value = 12000
$ws.Range("K136").Value = 6804
$ws.Range("L136").Value = 36000
$ws.Range("M136").Value = -4254
$ws.Range("N136").Value = -41100
# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3485.6333
$ws.Range("I132").Value = 4622.467
$ws.Range("J132").Value = 2348.8
$ws.Range("K132").Value = 13867.401
$ws.Range("L132").Value = 7046.400000000001
$ws.Range("M132").Value = -11337.401
$ws.Range("N132").Value = -12106.4
$ws.Range("H136").Value = 38392
$ws.Range("I136").Value = 43124
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 129372
$ws.Range("L136").Value = 10000
$ws.Range("M136").Value = -126822
$ws.Range("N136").Value = -35100
